# Update file with new data
# The header row (Sno, IP, Status, Link, last_analysis_stats, Country,
# whois_date, Last_Modification_Date, AS_Owner) is removed, the four data
# rows shift up by one, and each row's hyperlink (column D) now targets the
# shifted-up cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the hyperlink target URLs (in row order) before they move, so we
# can re-create them pointing at the correct cells after the header row is
# removed.
$linkUrls = @()
foreach ($h in $ws.Hyperlinks) {
  $linkUrls += $h.Address
}

# Remove the header row entirely; data rows 2-5 become rows 1-4.
$ws.Rows.Item(1).Delete()

# The old hyperlink objects still reference their pre-shift rows (D2:D5), so
# drop them and rebuild against the new row positions (D1:D4).
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("D1"), $linkUrls[0])
$ws.Hyperlinks.Add($ws.Range("D2"), $linkUrls[1])
$ws.Hyperlinks.Add($ws.Range("D3"), $linkUrls[2])
$ws.Hyperlinks.Add($ws.Range("D4"), $linkUrls[3])

# Adding a hyperlink re-applies its own ad-hoc font formatting; restore the
# original shared "Hyperlink" cell style so the four link cells keep using
# the same style as before.
$ws.Range("D1:D4").Style = "Hyperlink"
